$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 37 with new values
$ws.Range("A37").Value = 46011
$ws.Range("B37").Value = 2

# Add new rows 38-41 with data
$ws.Range("A38").Value = 46012
$ws.Range("B38").Value = 13

$ws.Range("A39").Value = 46010
$ws.Range("B39").Value = 2

$ws.Range("A40").Value = 46014
$ws.Range("B40").Value = 43

$ws.Range("A41").Value = 46013
$ws.Range("B41").Value = 1

# Row 42 gets the old row 37 data (date 46009, count 4), moved to the end
$ws.Range("A42").Value = 46009
$ws.Range("B42").Value = 4

# Copy style (date format) from A37 down to A38:A42
$ws.Range("A37").Copy()
$ws.Range("A38:A42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection / scroll position to match the final view
$ws.Range("A37:B42").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
